$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (D = Fecha, J = Volumen, K = Precio minimo,
# L = Precio maximo, M = Precio promedio ponderado, N = Unidad de
# comercializacion, O = Origen, P = Precio $/Kg, Q = Kg o Unidades).
# The rows are being rotated: row2<-row5, row3<-row4, row4<-row2, row5<-row3.

$rows = @{
    2 = @{ D = 44221; J = 250; K = 1300; L = 1500; M = 1420; N = "`$/atado"; O = "Provincia de Diguillín"; P = 1420; Q = 1 }
    3 = @{ D = 44687; J = 160; K = 3000; L = 3500; M = 3250; N = "`$/docena de matas"; O = "Región Metropolitana"; P = 542; Q = 6 }
    4 = @{ D = 44692; J = 120; K = 3000; L = 3500; M = 3250; N = "`$/docena de matas"; O = "Región Metropolitana"; P = 542; Q = 6 }
    5 = @{ D = 44691; J = 100; K = 3000; L = 3500; M = 3250; N = "`$/docena de matas"; O = "Región Metropolitana"; P = 542; Q = 6 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
}
